$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B2").Value = "https://interop.esante.gouv.fr/ig/fhir/tde/ValueSet/TypeCarteVS"
$wsMeta.Range("B3").Value = "2.0.0"
$wsMeta.Range("B8").Value = "2026-01-15T15:23:39+00:00"

$wsInclude = $wb.Worksheets.Item("Include #0")
$wsInclude.Range("B4").Value = "https://interop.esante.gouv.fr/ig/fhir/tde/CodeSystem/type-carte-code-system"
